$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'264.65"
$ws.Range("D3").Value = "'22.74"
$ws.Range("D4").Value = "'6.276"
$ws.Range("D5").Value = "'0.06156"
$ws.Range("D6").Value = "'3.593"
$ws.Range("D7").Value = "'6.715"
$ws.Range("D9").Value = "'0.8296"
$ws.Range("D11").Value = "'0.1601"
$ws.Range("D12").Value = "'0.08203"
$ws.Range("D14").Value = "'0.03139"
$ws.Range("D15").Value = "'0.09247"
$ws.Range("D16").Value = "'3.912"
$ws.Range("D17").Value = "'0.001704"
$ws.Range("D18").Value = "'0.04791"
$ws.Range("D19").Value = "'0.006281"
$ws.Range("D20").Value = "'0.005924"
$ws.Range("D21").Value = "'0.001105"
$ws.Range("D23").Value = "'3.766"
$ws.Range("D24").Value = "'2.304"
$ws.Range("D25").Value = "'0.3346"
$ws.Range("D27").Value = "'0.0002683"
$ws.Range("D40").Value = "'0.04641"
$ws.Range("D41").Value = "'0.006951"
$ws.Range("D42").Value = "'0.1135"
$ws.Range("D43").Value = "'0.003402"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("E44").Value = "43LocalTradersLCTWorstin24h"
$ws.Range("D45").Value = "'0.00006174"
$ws.Range("D47").Value = "'0.7784"
$ws.Range("D48").Value = "'0.2028"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.01241"
